$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 55.62696966666667
$ws.Cells.Item(2, 8).Value = 166.880909
$ws.Cells.Item(2, 9).Value = 0.2670320042914472
$ws.Cells.Item(2, 10).Value = 0.2670320042914472
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 7.602732666666667
$ws.Cells.Item(2, 14).Value = 22.808198
$ws.Cells.Item(2, 15).Value = 0.2697011854344982
$ws.Cells.Item(2, 16).Value = 0.2697011854344982
$ws.Cells.Item(2, 17).Value = 422.9169794324425
$ws.Cells.Item(2, 18).Value = 3806.252814891982
$ws.Cells.Item(2, 19).Value = 0.07201884810635331
$ws.Cells.Item(2, 20).Value = 0.07201884810635331

$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 55.62696966666667
$ws.Cells.Item(3, 8).Value = 166.880909
$ws.Cells.Item(3, 9).Value = 0.2670320042914472
$ws.Cells.Item(3, 10).Value = 0.2670320042914472
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 7.621805666666667
$ws.Cells.Item(3, 14).Value = 22.865417
$ws.Cells.Item(3, 15).Value = 0.2703777856696143
$ws.Cells.Item(3, 16).Value = 0.2703777856696144
$ws.Cells.Item(3, 17).Value = 423.9779526248948
$ws.Cells.Item(3, 18).Value = 3815.801573624053
$ws.Cells.Item(3, 19).Value = 0.07219952202324044
$ws.Cells.Item(3, 20).Value = 0.07219952202324047

$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 55.62696966666667
$ws.Cells.Item(4, 8).Value = 166.880909
$ws.Cells.Item(4, 9).Value = 0.2670320042914472
$ws.Cells.Item(4, 10).Value = 0.2670320042914472
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 10.36505866666667
$ws.Cells.Item(4, 14).Value = 31.095176
$ws.Cells.Item(4, 15).Value = 0.3676926089686856
$ws.Cells.Item(4, 16).Value = 0.3676926089686856
$ws.Cells.Item(4, 17).Value = 576.5768040438871
$ws.Cells.Item(4, 18).Value = 5189.191236394984
$ws.Cells.Item(4, 19).Value = 0.09818569433605945
$ws.Cells.Item(4, 20).Value = 0.09818569433605946

$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 55.62696966666667
$ws.Cells.Item(5, 8).Value = 166.880909
$ws.Cells.Item(5, 9).Value = 0.2670320042914472
$ws.Cells.Item(5, 10).Value = 0.2670320042914472
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 2.59987
$ws.Cells.Item(5, 14).Value = 7.79961
$ws.Cells.Item(5, 15).Value = 0.09222841992720189
$ws.Cells.Item(5, 16).Value = 0.0922284199272019
$ws.Cells.Item(5, 17).Value = 144.6228896272767
$ws.Cells.Item(5, 18).Value = 1301.60600664549
$ws.Cells.Item(5, 19).Value = 0.02462793982579397
$ws.Cells.Item(5, 20).Value = 0.02462793982579397

$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 63.357325
$ws.Cells.Item(6, 8).Value = 190.071975
$ws.Cells.Item(6, 9).Value = 0.3041408435993349
$ws.Cells.Item(6, 10).Value = 0.3041408435993349
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 7.602732666666667
$ws.Cells.Item(6, 14).Value = 22.808198
$ws.Cells.Item(6, 15).Value = 0.2697011854344982
$ws.Cells.Item(6, 16).Value = 0.2697011854344982
$ws.Cells.Item(6, 17).Value = 481.6888044501167
$ws.Cells.Item(6, 18).Value = 4335.199240051051
$ws.Cells.Item(6, 19).Value = 0.08202714605778894
$ws.Cells.Item(6, 20).Value = 0.08202714605778894

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 63.357325
$ws.Cells.Item(7, 8).Value = 190.071975
$ws.Cells.Item(7, 9).Value = 0.3041408435993349
$ws.Cells.Item(7, 10).Value = 0.3041408435993349
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 7.621805666666667
$ws.Cells.Item(7, 14).Value = 22.865417
$ws.Cells.Item(7, 15).Value = 0.2703777856696143
$ws.Cells.Item(7, 16).Value = 0.2703777856696144
$ws.Cells.Item(7, 17).Value = 482.8972187098417
$ws.Cells.Item(7, 18).Value = 4346.074968388575
$ws.Cells.Item(7, 19).Value = 0.08223292782407667
$ws.Cells.Item(7, 20).Value = 0.0822329278240767

$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 63.357325
$ws.Cells.Item(8, 8).Value = 190.071975
$ws.Cells.Item(8, 9).Value = 0.3041408435993349
$ws.Cells.Item(8, 10).Value = 0.3041408435993349
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 10.36505866666667
$ws.Cells.Item(8, 14).Value = 31.095176
$ws.Cells.Item(8, 15).Value = 0.3676926089686856
$ws.Cells.Item(8, 16).Value = 0.3676926089686856
$ws.Cells.Item(8, 17).Value = 656.7023905880667
$ws.Cells.Item(8, 18).Value = 5910.3215152926
$ws.Cells.Item(8, 19).Value = 0.1118303402769764
$ws.Cells.Item(8, 20).Value = 0.1118303402769764

$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 63.357325
$ws.Cells.Item(9, 8).Value = 190.071975
$ws.Cells.Item(9, 9).Value = 0.3041408435993349
$ws.Cells.Item(9, 10).Value = 0.3041408435993349
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 2.59987
$ws.Cells.Item(9, 14).Value = 7.79961
$ws.Cells.Item(9, 15).Value = 0.09222841992720189
$ws.Cells.Item(9, 16).Value = 0.0922284199272019
$ws.Cells.Item(9, 17).Value = 164.72080854775
$ws.Cells.Item(9, 18).Value = 1482.48727692975
$ws.Cells.Item(9, 19).Value = 0.02805042944049289
$ws.Cells.Item(9, 20).Value = 0.0280504294404929

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 52.65915966666667
$ws.Cells.Item(10, 8).Value = 157.977479
$ws.Cells.Item(10, 9).Value = 0.2527853132096735
$ws.Cells.Item(10, 10).Value = 0.2527853132096735
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 7.602732666666667
$ws.Cells.Item(10, 14).Value = 22.808198
$ws.Cells.Item(10, 15).Value = 0.2697011854344982
$ws.Cells.Item(10, 16).Value = 0.2697011854344982
$ws.Cells.Item(10, 17).Value = 400.3535133969825
$ws.Cells.Item(10, 18).Value = 3603.181620572842
$ws.Cells.Item(10, 19).Value = 0.06817649863307984
$ws.Cells.Item(10, 20).Value = 0.06817649863307984

$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 52.65915966666667
$ws.Cells.Item(11, 8).Value = 157.977479
$ws.Cells.Item(11, 9).Value = 0.2527853132096735
$ws.Cells.Item(11, 10).Value = 0.2527853132096735
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 7.621805666666667
$ws.Cells.Item(11, 14).Value = 22.865417
$ws.Cells.Item(11, 15).Value = 0.2703777856696143
$ws.Cells.Item(11, 16).Value = 0.2703777856696144
$ws.Cells.Item(11, 17).Value = 401.3578815493049
$ws.Cells.Item(11, 18).Value = 3612.220933943744
$ws.Cells.Item(11, 19).Value = 0.06834753323543143
$ws.Cells.Item(11, 20).Value = 0.06834753323543144

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 52.65915966666667
$ws.Cells.Item(12, 8).Value = 157.977479
$ws.Cells.Item(12, 9).Value = 0.2527853132096735
$ws.Cells.Item(12, 10).Value = 0.2527853132096735
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 10.36505866666667
$ws.Cells.Item(12, 14).Value = 31.095176
$ws.Cells.Item(12, 15).Value = 0.3676926089686856
$ws.Cells.Item(12, 16).Value = 0.3676926089686856
$ws.Cells.Item(12, 17).Value = 545.8152792823672
$ws.Cells.Item(12, 18).Value = 4912.337513541304
$ws.Cells.Item(12, 19).Value = 0.09294729132303116
$ws.Cells.Item(12, 20).Value = 0.09294729132303117

$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 52.65915966666667
$ws.Cells.Item(13, 8).Value = 157.977479
$ws.Cells.Item(13, 9).Value = 0.2527853132096735
$ws.Cells.Item(13, 10).Value = 0.2527853132096735
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 2.59987
$ws.Cells.Item(13, 14).Value = 7.79961
$ws.Cells.Item(13, 15).Value = 0.09222841992720189
$ws.Cells.Item(13, 16).Value = 0.0922284199272019
$ws.Cells.Item(13, 17).Value = 136.9069694425767
$ws.Cells.Item(13, 18).Value = 1232.16272498319
$ws.Cells.Item(13, 19).Value = 0.02331399001813102
$ws.Cells.Item(13, 20).Value = 0.02331399001813102

$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 36.672286
$ws.Cells.Item(14, 8).Value = 110.016858
$ws.Cells.Item(14, 9).Value = 0.1760418388995444
$ws.Cells.Item(14, 10).Value = 0.1760418388995444
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 7.602732666666667
$ws.Cells.Item(14, 14).Value = 22.808198
$ws.Cells.Item(14, 15).Value = 0.2697011854344982
$ws.Cells.Item(14, 16).Value = 0.2697011854344982
$ws.Cells.Item(14, 17).Value = 278.8095867335427
$ws.Cells.Item(14, 18).Value = 2509.286280601884
$ws.Cells.Item(14, 19).Value = 0.04747869263727609
$ws.Cells.Item(14, 20).Value = 0.04747869263727609

$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 36.672286
$ws.Cells.Item(15, 8).Value = 110.016858
$ws.Cells.Item(15, 9).Value = 0.1760418388995444
$ws.Cells.Item(15, 10).Value = 0.1760418388995444
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 7.621805666666667
$ws.Cells.Item(15, 14).Value = 22.865417
$ws.Cells.Item(15, 15).Value = 0.2703777856696143
$ws.Cells.Item(15, 16).Value = 0.2703777856696144
$ws.Cells.Item(15, 17).Value = 279.5090372444207
$ws.Cells.Item(15, 18).Value = 2515.581335199786
$ws.Cells.Item(15, 19).Value = 0.04759780258686581
$ws.Cells.Item(15, 20).Value = 0.04759780258686581

$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 36.672286
$ws.Cells.Item(16, 8).Value = 110.016858
$ws.Cells.Item(16, 9).Value = 0.1760418388995444
$ws.Cells.Item(16, 10).Value = 0.1760418388995444
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 10.36505866666667
$ws.Cells.Item(16, 14).Value = 31.095176
$ws.Cells.Item(16, 15).Value = 0.3676926089686856
$ws.Cells.Item(16, 16).Value = 0.3676926089686856
$ws.Cells.Item(16, 17).Value = 380.1103958307787
$ws.Cells.Item(16, 18).Value = 3420.993562477008
$ws.Cells.Item(16, 19).Value = 0.06472928303261853
$ws.Cells.Item(16, 20).Value = 0.06472928303261855

$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 36.672286
$ws.Cells.Item(17, 8).Value = 110.016858
$ws.Cells.Item(17, 9).Value = 0.1760418388995444
$ws.Cells.Item(17, 10).Value = 0.1760418388995444
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 2.59987
$ws.Cells.Item(17, 14).Value = 7.79961
$ws.Cells.Item(17, 15).Value = 0.09222841992720189
$ws.Cells.Item(17, 16).Value = 0.0922284199272019
$ws.Cells.Item(17, 17).Value = 95.34317620282
$ws.Cells.Item(17, 18).Value = 858.08858582538
$ws.Cells.Item(17, 19).Value = 0.01623606064278401
$ws.Cells.Item(17, 20).Value = 0.01623606064278401

